$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Karkosz et al. (2024)"
$ws.Range("B7").Value = 0.803826184793257
$ws.Range("C7").Value = 0.973169571781571

$ws.Range("A8").Value = "Sun et al. (2023)"
$ws.Range("B8").Value = 0.0916082071125613
$ws.Range("C8").Value = 0.999916682001268

$ws.Range("A9").Value = "Baez et al. (2017)"
$ws.Range("B9").Value = 0.826887932363078
$ws.Range("C9").Value = 0.453098578130092
